$d = $word.ActiveDocument

# Find the paragraph that starts with "Sono classi statiche" -- two new
# paragraphs (about Extract Class / Extract Method for the view refactor)
# are inserted immediately before it. The trailing "_GoBack" bookmark that
# used to sit on the very last paragraph ends up on the new second
# paragraph instead, since that is now the most-recently-edited spot.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Sono classi statiche*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphBefore()
$target.Range.InsertParagraphBefore()

$p1 = $d.Paragraphs.Item($targetIndex)
$p2 = $d.Paragraphs.Item($targetIndex + 1)
$p3 = $d.Paragraphs.Item($targetIndex + 2)
$p4 = $d.Paragraphs.Item($targetIndex + 3)

$p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Sono tutti </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Extract</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Class (del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>r</w:t></w:r><w:r><w:t>efactoring</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) ?</w:t></w:r></w:p>')
$p2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Anche per la grafica sono tutti </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>extract</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>method</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>syso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> diventano metodi nella view</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
$p3.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Sono classi statiche, cioè non vengono istanziate ma hanno metodi statici. </w:t></w:r></w:p>')
$p4.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Potrebbero diventare singleton? Così si possono istanziare ma solo una volta.</w:t></w:r></w:p>')
